# lab6.pptx - "Approach:" slide (slide 6) - update the light-sensor threshold
# from 200 to 250 in the sample code TextBox.
#
#   if (light < 200) {   ->   if (light < 250) {
#    delay(200);         ->    delay(250);

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(6)
$sh = $s.Shapes.Item("TextBox 6")
$tr = $sh.TextFrame.TextRange

# The textbox auto-fits its height to the text (<a:spAutoFit/>); remember
# the current height so it can be restored after the edits below, since
# the source edit only touched the code text, not the shape's geometry.
$origHeight = $sh.Height

# --- "if (light < 200) {" -> "if (light < 250) {" -----------------------
# This paragraph holds a single run; Replace() splits it into three runs
# around the matched substring, which is exactly how the slide author's
# edit was recorded (the surrounding "if (light " / ") {" text keeps the
# original run formatting, the replaced "< 250" becomes its own run).
$ifParagraph = $tr.Paragraphs(3, 1)
$ifParagraph.Replace("< 200", "< 250", 1, $false, $false) | Out-Null

# --- " delay(200);" -> " delay(250);" ------------------------------------
# This paragraph is a single run too; just rewrite its text in place so no
# extra run is introduced, matching the recorded edit.
$delayParagraph = $tr.Paragraphs(10, 1)
$delayRun = $delayParagraph.Runs(1, 1)
$delayRun.Text = " delay(250);"

# Restore the auto-fit height that the text edits recalculated.
$sh.Height = $origHeight
